# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-11 11:13:06
#
# Applies the following changes to "Session Analysis Results":
#   1. Swap the "Recorded By" text order from "dnasr281@gmail.com, System"
#      to "System, dnasr281@gmail.com" for the rows whose session has
#      actually been recorded.
#   2. Update the "Missing Sessions" / "Pending Sessions" statistics
#      (L7, L8) and the per-group P/Q attendance counts that shifted as a
#      result of newly recorded sessions.
#   3. Move the six still-outstanding "B1-*" sessions whose date
#      (11/01/2026) has now passed from "Pending" to "Not Recorded",
#      which also flips their row shading from yellow to pink (the style
#      already used elsewhere in the sheet for "Not Recorded" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. "Recorded By" column (G) text swap
# ---------------------------------------------------------------------
$recordedByRows = @(2,3,4,5,6,7,16,17,22,23,37,38,43,44,58,59,64,65,79,80,85,86,87,88,89,90,
  99,100,105,106,107,108,109,110,119,120,125,126,127,128,129,130,139,140,
  145,146,147,148,149,150,159,160,165,166,167,168,169,170,179,180,185,186,
  200,201,206,207,221,222,227,228,242,243)

foreach ($r in $recordedByRows) {
  $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------
# 2. Statistics updates
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 72
$ws.Range("L8").Value = 54

$pqRows = @(16,17,18,24,25,26)
foreach ($r in $pqRows) {
  $ws.Range("P$r").Value = 6
  $ws.Range("Q$r").Value = 4
}

# ---------------------------------------------------------------------
# 3. "Pending" -> "Not Recorded" for sessions dated 11/01/2026
#    (copy the formatting already used for "Not Recorded" rows, e.g.
#    row 10, so the fill/style index is reused rather than duplicated)
# ---------------------------------------------------------------------
$notRecordedRows = @(32,53,74,195,216,237)
foreach ($r in $notRecordedRows) {
  $ws.Range("A10:I10").Copy()
  $ws.Range("A${r}:I${r}").PasteSpecial(-4122)
  $ws.Range("I$r").Value = "Not Recorded"
}
